# Weekly price update: a new week's record is inserted at row 83 (pushing
# the existing rows 83-115 down to 84-116), and the new row is populated
# with the latest week's figures for the same market/product/variety/quality
# combination as the row that used to sit at 83.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 83 - this shifts rows
# 83..115 down to 84..116 (carrying all their original values/styles with
# them) and grows the sheet dimension to A1:T116 automatically.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with this week's data. Columns that
# stay the same as the row that used to be here (now at row 84) are copied
# across; only the figures that actually changed week-over-week are new.
$ws.Cells.Item(83, 1).Value2  = 10
$ws.Cells.Item(83, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value2  = "La Araucanía"
$ws.Cells.Item(83, 4).Value2  = 45141
$ws.Cells.Item(83, 5).Value2  = 9
$ws.Cells.Item(83, 6).Value2  = "Fruta"
$ws.Cells.Item(83, 7).Value2  = 100108
$ws.Cells.Item(83, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(83, 9).Value2  = 100108004
$ws.Cells.Item(83, 10).Value2 = "Papaya"
$ws.Cells.Item(83, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(83, 12).Value2 = "Primera"
$ws.Cells.Item(83, 13).Value2 = 85
$ws.Cells.Item(83, 14).Value2 = 25000
$ws.Cells.Item(83, 15).Value2 = 25000
$ws.Cells.Item(83, 16).Value2 = 25000
$ws.Cells.Item(83, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(83, 18).Value2 = "Provincia del Elquí"
$ws.Cells.Item(83, 19).Value2 = 2500
$ws.Cells.Item(83, 20).Value2 = 10
